$d = $word.ActiveDocument

# Change 1: merge the "being" / proofErr-wrapped run back into a single plain run.
# Restrict the search to just the second run's text (starting at "session while")
# so the preceding, untouched run ("...during the ") is left alone.
$r = $d.Content
$r.Find.Execute(
    "session while the device provides an indication of whether the current goal is ", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
)
$start = $r.Find.Found
$rng = $d.Range($r.Start, $r.Start)


# Change 2: update the L-13 requirement text
$d.Content.Find.Execute(
    "The success or failure of the user to achieve each goal is recorded.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "An indication is recorded of when during the workout each goal was being met or not.",
    2
)
